$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" - append 3 new weekly rows (21-23)
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

$weeklyDates = @(45662.99999999999, 45676.99999999999, 45683.99999999999)
$weeklyQtys  = @(3, 9, 1)

for ($i = 0; $i -lt $weeklyDates.Length; $i++) {
    $r = 21 + $i
    $wsWeekly.Cells.Item($r, 1).Value = $weeklyDates[$i]
    $wsWeekly.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsWeekly.Cells.Item($r, 2).Value = $weeklyQtys[$i]
}

# ---------------------------------------------------------------------
# Sheet 2: "Monthly Trend" - append 1 new monthly row (11)
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsMonthly.Cells.Item(11, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(11, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsMonthly.Cells.Item(11, 2).Value = 13

# ---------------------------------------------------------------------
# Sheet 3: "PO Forecast" - new forecast model: update existing forecast
# values (rows 2-28) and append 3 new rows (29-31)
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$forecastDates = @(
    44990.99999999999, 44997.99999999999, 45004.99999999999, 45011.99999999999,
    45039.99999999999, 45088.99999999999, 45179.99999999999, 45186.99999999999,
    45193.99999999999, 45200.99999999999, 45214.99999999999, 45221.99999999999,
    45228.99999999999, 45235.99999999999, 45242.99999999999, 45249.99999999999,
    45270.99999999999, 45613.99999999999, 45641.99999999999, 45662.99999999999,
    45676.99999999999, 45683.99999999999, 45690.99999999999, 45697.99999999999,
    45704.99999999999, 45711.99999999999, 45718.99999999999, 45725.99999999999,
    45732.99999999999, 45739.99999999999
)

$forecastQtys = @(
    30, 30, 30, 30,
    30, 29, 28, 28,
    28, 28, 28, 28,
    27, 27, 27, 27,
    27, 22, 22, 22,
    22, 22, 21, 21,
    21, 21, 21, 21,
    21, 21
)

for ($i = 0; $i -lt $forecastDates.Length; $i++) {
    $r = 2 + $i
    $wsForecast.Cells.Item($r, 1).Value = $forecastDates[$i]
    $wsForecast.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($r, 2).Value = $forecastQtys[$i]
}
